$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-03-31"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 03-31)"

# Update March total (row 4) and grand Total row (row 14) for the "Total" column (I)
$ws.Range("I4").Value = 133
$ws.Range("I14").Value = 433
